$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Row 88: new "FRIDAY" section header row (copy exact formatting+value from row 41) ---
$ws.Range("A41:F41").Copy($ws.Range("A88:F88"))

# --- Row 89: Pickup PC / WC / 001-DH ---
$ws.Range("A76:F76").Copy($ws.Range("A89:F89"))
$ws.Range("A89").Value = "Pickup PC"
$ws.Range("B89").Value = 42608
$ws.Range("C89").Value = "1530"
$ws.Range("D89").Value = "WC"
$ws.Range("E89").Value = "001-DH"
$ws.Range("F89").Value = "Pick up PC and Projector cart and large portable screen also. PICK UP MATS. Return equipment to Vanier 040 storeroom."

# --- Row 90: Pickup Large PA / WC / 001-DH ---
$ws.Range("A76:F76").Copy($ws.Range("A90:F90"))
$ws.Range("A90").Value = "Pickup Large PA"
$ws.Range("B90").Value = 42608
$ws.Range("C90").Value = "1530"
$ws.Range("D90").Value = "WC"
$ws.Range("E90").Value = "001-DH"
$ws.Range("F90").Value = "Lectern mic, stand and cables with large PA system. To Vanier 040."
$ws.Rows.Item(90).RowHeight = 15

# --- Row 91: Pickup PC / VC / 001-JCR ---
$ws.Range("A76:F76").Copy($ws.Range("A91:F91"))
$ws.Range("A91").Value = "Pickup PC"
$ws.Range("B91").Value = 42608
$ws.Range("C91").Value = "1800"
$ws.Range("D91").Value = "VC"
$ws.Range("E91").Value = "001-JCR"
$ws.Range("F91").Value = "May include portable screen. All to Vanier 040."
$ws.Rows.Item(91).RowHeight = 15

# --- Row 92: Pickup Small PA / VC / 001-JCR ---
$ws.Range("A76:F76").Copy($ws.Range("A92:F92"))
$ws.Range("A92").Value = "Pickup Small PA"
$ws.Range("B92").Value = 42608
$ws.Range("C92").Value = "1800"
$ws.Range("D92").Value = "VC"
$ws.Range("E92").Value = "001-JCR"
$ws.Range("F92").Value = "Lectern mic with small PA etc. Return to Vanier 040."

# --- Row 93: Other / WC / 283B (same task as row 76/87, next day) ---
$ws.Range("A76:F76").Copy($ws.Range("A93:F93"))
$ws.Range("B93").Value = 42608

# --- Update the frozen-pane view state to match the author's final scroll/selection ---
$ws.Application.GoTo($ws.Range("A75"))
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E101:E102").Select()
